$wb = $excel.ActiveWorkbook

# --- Sheet: 2o Parcial ---
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Cells.Item(8, 5).Value = 18
$ws2.Cells.Item(8, 6).Value = 6
$ws2.Cells.Item(8, 7).Value = 75
$ws2.Cells.Item(8, 8).Value = 25
$ws2.Cells.Item(8, 9).Value = 7.8
$ws2.Cells.Item(8, 10).Value = 0
$ws2.Cells.Item(8, 11).Value = 0

$ws2.Cells.Item(9, 5).Value = 32
$ws2.Cells.Item(9, 6).Value = 2
$ws2.Cells.Item(9, 7).Value = 94.1
$ws2.Cells.Item(9, 8).Value = 5.9
$ws2.Cells.Item(9, 9).Value = 8.5
$ws2.Cells.Item(9, 10).Value = 0
$ws2.Cells.Item(9, 11).Value = 0

$ws2.Cells.Item(10, 5).Value = 31
$ws2.Cells.Item(10, 6).Value = 0
$ws2.Cells.Item(10, 7).Value = 100
$ws2.Cells.Item(10, 8).Value = 0
$ws2.Cells.Item(10, 9).Value = 9
$ws2.Cells.Item(10, 10).Value = 0
$ws2.Cells.Item(10, 11).Value = 0

$ws2.Cells.Item(11, 5).Value = 40
$ws2.Cells.Item(11, 6).Value = 0
$ws2.Cells.Item(11, 7).Value = 100
$ws2.Cells.Item(11, 8).Value = 0
$ws2.Cells.Item(11, 9).Value = 9.3
$ws2.Cells.Item(11, 10).Value = 0
$ws2.Cells.Item(11, 11).Value = 0

$ws2.Cells.Item(12, 5).Value = 19
$ws2.Cells.Item(12, 6).Value = 5
$ws2.Cells.Item(12, 7).Value = 79.2
$ws2.Cells.Item(12, 8).Value = 20.8
$ws2.Cells.Item(12, 9).Value = 7.6
$ws2.Cells.Item(12, 10).Value = 0
$ws2.Cells.Item(12, 11).Value = 0

$ws2.Cells.Item(13, 5).Value = 35
$ws2.Cells.Item(13, 6).Value = 3
$ws2.Cells.Item(13, 7).Value = 92.1
$ws2.Cells.Item(13, 8).Value = 7.9
$ws2.Cells.Item(13, 9).Value = 8.9
$ws2.Cells.Item(13, 10).Value = 0
$ws2.Cells.Item(13, 11).Value = 0

$ws2.Cells.Item(14, 5).Value = 175
$ws2.Cells.Item(14, 6).Value = 16
$ws2.Cells.Item(14, 7).Value = 91.6
$ws2.Cells.Item(14, 8).Value = 8.4
$ws2.Cells.Item(14, 9).Value = 8.5
$ws2.Cells.Item(14, 10).Value = 0
$ws2.Cells.Item(14, 11).Value = 0

$ws2.Cells.Item(15, 5).Value = 233
$ws2.Cells.Item(15, 6).Value = 39
$ws2.Cells.Item(15, 7).Value = 85.7
$ws2.Cells.Item(15, 8).Value = 14.3
$ws2.Cells.Item(15, 9).Value = 7.7
$ws2.Cells.Item(15, 10).Value = 0
$ws2.Cells.Item(15, 11).Value = 0

# --- Sheet: Final ---
$ws3 = $wb.Worksheets.Item("Final")
$ws3.Cells.Item(8, 9).Value = 7.3

$ws3.Cells.Item(9, 9).Value = 8.4

$ws3.Cells.Item(10, 5).Value = 31
$ws3.Cells.Item(10, 6).Value = 0
$ws3.Cells.Item(10, 7).Value = 100
$ws3.Cells.Item(10, 8).Value = 0
$ws3.Cells.Item(10, 9).Value = 9

$ws3.Cells.Item(11, 9).Value = 9.3

$ws3.Cells.Item(12, 5).Value = 19
$ws3.Cells.Item(12, 6).Value = 5
$ws3.Cells.Item(12, 7).Value = 79.2
$ws3.Cells.Item(12, 8).Value = 20.8
$ws3.Cells.Item(12, 9).Value = 7.5

$ws3.Cells.Item(13, 5).Value = 35
$ws3.Cells.Item(13, 6).Value = 3
$ws3.Cells.Item(13, 7).Value = 92.1
$ws3.Cells.Item(13, 8).Value = 7.9
$ws3.Cells.Item(13, 9).Value = 9

$ws3.Cells.Item(14, 5).Value = 175
$ws3.Cells.Item(14, 6).Value = 16
$ws3.Cells.Item(14, 7).Value = 91.6
$ws3.Cells.Item(14, 8).Value = 8.4
$ws3.Cells.Item(14, 9).Value = 8.4

$ws3.Cells.Item(15, 5).Value = 233
$ws3.Cells.Item(15, 6).Value = 39
$ws3.Cells.Item(15, 7).Value = 85.7
$ws3.Cells.Item(15, 8).Value = 14.3
$ws3.Cells.Item(15, 9).Value = 7.5
